# Add "shipheader equip part 1" rows to the condition sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28: itemNotNil condition
$ws.Range("A28").Value = "itemNotNil"
$ws.Range("B28").Value = "选择了装备"
$ws.Range("C28").Value = "data"
$ws.Range("D28").Value = "item"
$ws.Range("E28").Value = "!="
$ws.Range("F28").Value = "nil"

# Row 29: shipHasHeader condition
$ws.Range("A29").Value = "shipHasHeader"
$ws.Range("B29").Value = "装备了船首像"
$ws.Range("C29").Value = "data"
$ws.Range("D29").Value = "ship"
$ws.Range("E29").Value = "shipHeader"
$ws.Range("F29").Value = ";"
$ws.Range("G29").Value = ";"

# Update the active selection to match the committed workbook state.
$ws.Range("E27").Select() | Out-Null
